# Add two new duck rows (Strawberry, Unicorn Large) to the Ducks table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducks")
$lo = $ws.ListObjects.Item(1)

# --- Row 1: Strawberry ---
$row1 = $lo.ListRows.Add()
$r1 = $row1.Range.Row

# --- Row 2: Unicorn Large ---
$row2 = $lo.ListRows.Add()
$r2 = $row2.Range.Row

# Fill column A (Duck name) for both new rows first, so the shared-string
# table gets "Strawberry" / "Unicorn Large" allocated before "Target".
$ws.Cells.Item($r1, 1).Value = "Strawberry"
$ws.Cells.Item($r2, 1).Value = "Unicorn Large"

# --- Row 1 remaining fields ---
$ws.Cells.Item($r1, 3).Value = "Phyiscal Store"
$ws.Cells.Item($r1, 4).Value = "Target"
$ws.Cells.Item($r1, 5).Value = "Burlington"
$ws.Cells.Item($r1, 6).Value = "MA"
$ws.Cells.Item($r1, 7).Value = "USA"
$ws.Cells.Item($r1, 8).Value = "USA"

$ws.Range("I40").Copy()
$ws.Cells.Item($r1, 9).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item($r1, 9).Value2 = 44933

$ws.Range("K2").Copy()
$ws.Range($ws.Cells.Item($r1, 10), $ws.Cells.Item($r1, 11)).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item($r1, 10).Value2 = 42.483620999999999
$ws.Cells.Item($r1, 11).Value2 = -71.185913999999997
$ws.Cells.Item($r1, 12).WrapText = $true
$ws.Cells.Item($r1, 13).Value = "Julia"
$ws.Cells.Item($r1, 14).Value2 = 1
$ws.Cells.Item($r1, 15).Value2 = 44
$ws.Cells.Item($r1, 16).Value2 = 8.4
$ws.Cells.Item($r1, 17).Value2 = 7.4
$ws.Cells.Item($r1, 18).Value2 = 7.6

# --- Row 2 remaining fields ---
$ws.Cells.Item($r2, 3).Value = "Phyiscal Store"
$ws.Cells.Item($r2, 4).Value = "Target"
$ws.Cells.Item($r2, 5).Value = "Burlington"
$ws.Cells.Item($r2, 6).Value = "MA"
$ws.Cells.Item($r2, 7).Value = "USA"
$ws.Cells.Item($r2, 8).Value = "USA"

$ws.Range("I40").Copy()
$ws.Cells.Item($r2, 9).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item($r2, 9).Value2 = 44933

$ws.Range("K2").Copy()
$ws.Range($ws.Cells.Item($r2, 10), $ws.Cells.Item($r2, 11)).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item($r2, 10).Value2 = 42.483620999999999
$ws.Cells.Item($r2, 11).Value2 = -71.185913999999997
$ws.Cells.Item($r2, 12).WrapText = $true
$ws.Cells.Item($r2, 13).Value = "Julia"
$ws.Cells.Item($r2, 14).Value2 = 1
$ws.Cells.Item($r2, 15).Value2 = 46
$ws.Cells.Item($r2, 16).Value2 = 8.7
$ws.Cells.Item($r2, 17).Value2 = 7.3
$ws.Cells.Item($r2, 18).Value2 = 9

# Match the final on-screen selection/scroll state left by the edit.
$ws.Range("L43").Select()
